$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the LiveSLR copyright/build string in B2 (shared string content change)
$ws.Range("B2").Value = "Copyright @ 2022 Cytel Inc. LiveSLR 4.0.0.0 - Build #49237"

# Move the active selection from B2 to C2
$ws.Range("C2").Select()
